$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Playoff")

# Final match score (row 21): Skor 1 = 0, Skor 2 = 2
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 2

# 3.luk (third place) match score (row 24): Skor 1 = 4, Skor 2 = 2
$ws.Range("D24").Value = 4
$ws.Range("E24").Value = 2

# Leave selection on F21, matching the final cursor position recorded in the file
$ws.Range("F21").Select()
